# Apply crypto price/volume updates from the Sep 7 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of updated prices are purely numeric-looking strings (e.g. "215.96").
# Force those specific cells to remain Text so Excel does not auto-convert them
# to Number cells -- the source data stores every Price/Volume value as text.
$textCells = @("D5", "D10", "D17", "D21", "D22", "D23", "D26", "D27", "D28", "D29", "D31", "D38", "D39", "D43", "D47", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '25.843.52'
$ws.Range('E2').Value = '  -0.02%  '

$ws.Range('D3').Value = '1.640.73'
$ws.Range('E3').Value = '  +0.28%  '

$ws.Range('E4').Value = '  -0.16%  '

$ws.Range('D5').Value = '215.96'
$ws.Range('E5').Value = '  +0.31%  '

$ws.Range('E7').Value = '  -0.12%  '

$ws.Range('E8').Value = '  -0.26%  '

$ws.Range('E9').Value = '  -0.88%  '

$ws.Range('D10').Value = '19.74'
$ws.Range('E10').Value = '  -1.96%  '

$ws.Range('E11').Value = '  +1.69%  '

$ws.Range('E12').Value = '  +0.29%  '

$ws.Range('D13').Value = '1.867.14'
$ws.Range('E13').Value = '  +0.27%  '

$ws.Range('D14').Value = '1.639.13'
$ws.Range('E14').Value = '  -0.33%  '

$ws.Range('E15').Value = '  -0.14%  '

$ws.Range('E16').Value = '  +0.42%  '

$ws.Range('D17').Value = '63.15'
$ws.Range('E17').Value = '  -0.15%  '

$ws.Range('D18').Value = '25.874.43'
$ws.Range('E18').Value = '  +0.06%  '

$ws.Range('E19').Value = '  -0.14%  '

$ws.Range('E20').Value = '  +2.25%  '

$ws.Range('D21').Value = '193.14'
$ws.Range('E21').Value = '  -0.45%  '

$ws.Range('D22').Value = '10.00'
$ws.Range('E22').Value = '  +0.51%  '

$ws.Range('D23').Value = '6.36'
$ws.Range('E23').Value = '  +2.70%  '

$ws.Range('E24').Value = '  +4.46%  '

$ws.Range('D26').Value = '142.07'
$ws.Range('E26').Value = '  +2.61%  '

$ws.Range('D27').Value = '0.123'
$ws.Range('E27').Value = '  -0.32%  '

$ws.Range('D28').Value = '6.96'
$ws.Range('E28').Value = '  +1.83%  '

$ws.Range('D29').Value = '15.56'
$ws.Range('E29').Value = '  -0.04%  '

$ws.Range('E30').Value = '  -0.27%  '

$ws.Range('D31').Value = '0.0496'
$ws.Range('E31').Value = '  +0.28%  '

$ws.Range('E32').Value = '  +1.14%  '

$ws.Range('E33').Value = '  +0.51%  '

$ws.Range('E34').Value = '  +0.63%  '

$ws.Range('E35').Value = '  -0.27%  '

$ws.Range('E36').Value = '  +0.60%  '

$ws.Range('D37').Value = '1.132.10'
$ws.Range('E37').Value = '  +0.71%  '

$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '2.53'
$ws.Range('E38').Value = '  -2.24%  '

$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').Value = '0.548'
$ws.Range('E39').Value = '  -0.41%  '

$ws.Range('E40').Value = '  -0.48%  '

$ws.Range('E41').Value = '  +0.12%  '

$ws.Range('E42').Value = '  +1.25%  '

$ws.Range('D43').Value = '100.85'
$ws.Range('E43').Value = '  +1.43%  '

$ws.Range('E44').Value = '  +0.62%  '

$ws.Range('D45').Value = '1.776.29'
$ws.Range('E45').Value = '  +0.16%  '

$ws.Range('E46').Value = '  +4.37%  '

$ws.Range('D47').Value = '55.50'
$ws.Range('E47').Value = '  +0.02%  '

$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '0.419'
$ws.Range('E48').Value = '  -1.20%  '

$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '1.44'
$ws.Range('E49').Value = '  +6.22%  '

$ws.Range('E50').Value = '  -0.09%  '

$ws.Range('E51').Value = '  +3.22%  '
